# Weekly data refresh: insert the newest "Achicoria" price record as a new
# row right after the current row 25, pushing the existing rows 26-30 down
# to rows 27-31 (their contents are unchanged, just relocated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26; everything below shifts down.
$ws.Rows.Item(26).Insert()

# Populate the new row with the latest weekly observation.
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44932
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = 100112010
$ws.Range("G26").Value = "Achicoria"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 6500
$ws.Range("N26").Value = "$/caja 16 unidades"
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 406
$ws.Range("Q26").Value = 16
$ws.Range("R26").Value = "Hortaliza"

# Keep the new "Fecha" cell formatted the same way as the rest of column D.
$ws.Range("D26").NumberFormat = $ws.Range("D27").NumberFormat
